# Add a new row for the paper "Seasonality and Sex-Biased Fluctuation of Birth
# Weight in Tibetan Populations" before the existing row 82, shifting all
# subsequent rows down by one (old row 82 becomes new row 83, etc.).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new, blank row at position 82; everything below shifts down.
$ws.Rows.Item(82).Insert()

# Columns A, B, C, E, F, J, K, L, M are stored as text in this sheet (titles,
# URLs, type labels, dates-as-text, author names, and volume/issue numbers
# that look numeric but are text). Force text formatting before assigning so
# Excel does not auto-convert things like "2022-01-31" into a date serial or
# "2" / "1" into plain numbers.
$textCols = @("A", "B", "C", "E", "F", "J", "K", "L", "M")
foreach ($col in $textCols) {
    $ws.Range($col + "82").NumberFormat = "@"
}

$ws.Range("A82").Value = "Seasonality and Sex-Biased Fluctuation of Birth Weight in Tibetan Populations"
$ws.Range("B82").Value = "http://link.springer.com/article/10.1007/s43657-021-00038-7"
$ws.Range("C82").Value = "Article"
$ws.Range("D82").Value = 2022
$ws.Range("E82").Value = "2022-01"
$ws.Range("F82").Value = "2022-01-31"
$ws.Range("G82").Value = 1219
$ws.Range("H82").Value = 1
$ws.Range("I82").Value = 2
$ws.Range("J82").Value = " Ouzhuluobu, Xuebin Qi, Bing Su"
$ws.Range("K82").Value = "2"
$ws.Range("L82").Value = "1"
$ws.Range("M82").Value = "2023-09-18"
